# Update the agenda schedule on the active sheet to reflect the new
# topic ordering / durations (Break moved, new Topic 7 added, Duration
# header renamed, times shifted, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: E1 "60" -> "Duration"
$ws.Range("E1").Value = "Duration"

# Data rows: Start-time (B), End-time (C), Topic (D), Duration (E)
$rows = @(
    @{ Row = 2;  B = "10:00:00"; C = "10:25:00"; D = "Topic 1"; E = "25" },
    @{ Row = 3;  B = "10:25:00"; C = "10:55:00"; D = "Topic 2"; E = "30" },
    @{ Row = 4;  B = "10:55:00"; C = "11:25:00"; D = "Break";   E = "30" },
    @{ Row = 5;  B = "11:25:00"; C = "11:45:00"; D = "Topic 3"; E = "20" },
    @{ Row = 6;  B = "11:45:00"; C = "12:45:00"; D = "Topic 4"; E = "60" },
    @{ Row = 7;  B = "12:45:00"; C = "13:00:00"; D = "Break";   E = "15" },
    @{ Row = 8;  B = "13:00:00"; C = "13:25:00"; D = "Topic 5"; E = "25" },
    @{ Row = 9;  B = "13:25:00"; C = "13:55:00"; D = "Topic 6"; E = "30" },
    @{ Row = 10; B = "13:55:00"; C = "14:40:00"; D = "Lunch";   E = "45" },
    @{ Row = 11; B = "14:40:00"; C = "15:00:00"; D = "Topic 7"; E = "20" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D

    # Column E holds purely-numeric-looking text ("25", "30", ...). The
    # source workbook stores these as inline *strings*, not numbers, so
    # force text formatting before assigning, then drop back to the
    # Normal style so no stray per-cell formatting is left behind.
    $cell = $ws.Cells.Item($r.Row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $r.E
    $cell.Style = "Normal"
}
